$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header for intervention_type (copy header formatting from J1)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K1").Value = "intervention_type"

# Fill in intervention_type values for each row
$values = @("BEHAVIORAL", "BEHAVIORAL", "BEHAVIORAL", "PROCEDURE", "BEHAVIORAL", "BEHAVIORAL", "BEHAVIORAL", "PROCEDURE", "DIETARY_SUPPLEMENT", "OTHER")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
